# database : add deficiency files
# Insert 5 new quarterly columns (D:H) before the existing data, shifting the
# existing quarters from D:H -> I:M, then populate the new columns with the
# older quarters' figures (1399/03 .. 1400/03) that were missing before.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert 5 blank columns at D, pushing old D:H -> I:M
$ws.Columns("D:H").Insert()

# 2) Column widths for the freshly inserted D:H (mirrors the I:M pattern:
#    29,29,29 / 31 / 29 in "raw" xlsx width units -> ColumnWidth = raw - 0.83)
$ws.Columns("D:F").ColumnWidth = 28.17
$ws.Columns("G:G").ColumnWidth = 30.17
$ws.Columns("H:H").ColumnWidth = 28.17

# 3) Header row 8 (quarter labels) for the new columns
$ws.Range("D8").Value = "فصل اول منتهی به 1399/03"
$ws.Range("E8").Value = "فصل دوم منتهی به 1399/06"
$ws.Range("F8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("G8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("H8").Value = "فصل اول منتهی به 1400/03"

# 4) Row 9 (publish-date labels) for the new columns
$ws.Range("D9").Value = "1400-04-30 (3)"
$ws.Range("E9").Value = "1400-11-24 (4)"
$ws.Range("F9").Value = "1400-10-29 (2)"
$ws.Range("G9").Value = "1401-05-09 (9)"
$ws.Range("H9").Value = "1401-04-30 (2)"

# 5) Financial data for rows 11-27, columns D:H (new quarters)
$data = @{
    11 = @(126685284,146932639,220934628,279484194,336653925)
    12 = @(-68737017,-82274038,-112123208,-143168055,-164051934)
    13 = @(57948267,64658601,108811420,136316139,172601991)
    14 = @(-3830892,-3429424,-5062177,-6356640,-6317127)
    15 = @(0,0,0,0,0)
    16 = @(-891305,4761111,106691,2379648,-473711)
    17 = @(53226070,65990288,103855934,132339147,165811153)
    18 = @(-7048609,-4898949,-5822180,-12473570,-7785869)
    19 = @(381037,15740001,8301347,65587512,10915093)
    20 = @(46558498,76831340,106335101,185453089,168940377)
    21 = @(-7449360,-8987315,-6211230,-9755112,-30341156)
    22 = @(39109138,67844025,100123871,175697977,138599221)
    23 = @(0,0,0,0,0)
    24 = @(39109138,67844025,100123871,175697977,138599221)
    25 = @(133,232,342,600,473)
    26 = @(293000000,293000000,293000000,293000000,293000000)
    27 = @(74,128,189,332,262)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $col = 4 + $i
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

# 6) Row 26 ("سرمایه") previously used the "Comma" cell style (with thousands
#    separators) on its last two quarters and held stale 530,000,000 capital
#    figures; normalise the whole row to the plain right-aligned style used
#    by the rest of the row and the correct 293,000,000 capital value.
$ws.Range("J26").Copy()
$ws.Range("K26:M26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("L26").Value = 293000000
$ws.Range("M26").Value = 293000000
